# Updated cryptos list (price/volume refresh) — GitHub Actions data pull.
# All values are assigned as text (leading "'" forces text for pure
# numeric-looking strings) to preserve formatting such as trailing zeros
# and the existing inlineStr/text cell type used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.359.96"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("D3").Value = "1.842.03"
$ws.Range("E3").Value = "'  -0.27%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("D5").Value = "'239.20"
$ws.Range("E5").Value = "'  -0.53%  "
$ws.Range("D6").Value = "'0.6292"
$ws.Range("E6").Value = "'  -0.62%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E8").Value = "'  -0.51%  "
$ws.Range("D9").Value = "'0.2931"
$ws.Range("E9").Value = "'  -1.17%  "
$ws.Range("D10").Value = "'24.39"
$ws.Range("E10").Value = "'  -0.62%  "
$ws.Range("D11").Value = "'0.07702"
$ws.Range("E11").Value = "'  -0.16%  "
$ws.Range("D12").Value = "1.868.57"
$ws.Range("E12").Value = "'  -5.87%  "
$ws.Range("D13").Value = "'4.986"
$ws.Range("E13").Value = "'  +0.07%  "
$ws.Range("D14").Value = "'0.6779"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("D15").Value = "'0.00001040"
$ws.Range("E15").Value = "'  +4.41%  "
$ws.Range("D16").Value = "'82.81"
$ws.Range("E16").Value = "'  +0.07%  "
$ws.Range("D17").Value = "2.109.58"
$ws.Range("E17").Value = "'  -6.83%  "
$ws.Range("D18").Value = "'6.110"
$ws.Range("E18").Value = "'  -1.23%  "
$ws.Range("D19").Value = "29.387.51"
$ws.Range("E19").Value = "'  -0.10%  "
$ws.Range("D20").Value = "'227.42"
$ws.Range("E20").Value = "'  -1.72%  "
$ws.Range("D21").Value = "'12.42"
$ws.Range("E21").Value = "'  -0.71%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'7.421"
$ws.Range("E23").Value = "'  -2.16%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.07%  "
$ws.Range("D25").Value = "'156.56"
$ws.Range("E25").Value = "'  +1.17%  "
$ws.Range("E26").Value = "'  -0.57%  "
$ws.Range("D27").Value = "'8.349"
$ws.Range("E27").Value = "'  -0.90%  "
$ws.Range("E28").Value = "'  -0.46%  "
$ws.Range("D29").Value = "'1.455"
$ws.Range("E29").Value = "'  -1.05%  "
$ws.Range("D30").Value = "'1.275"
$ws.Range("E30").Value = "'  +1.24%  "
$ws.Range("D31").Value = "'0.05634"
$ws.Range("E31").Value = "'  -3.00%  "
$ws.Range("D32").Value = "'4.095"
$ws.Range("E33").Value = "'  -0.01%  "
$ws.Range("D34").Value = "'1.826"
$ws.Range("E34").Value = "'  -2.15%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("E35").Value = "'  -0.35%  "
$ws.Range("D36").Value = "'0.7061"
$ws.Range("E36").Value = "'  -1.44%  "
$ws.Range("D37").Value = "'2.588"
$ws.Range("E37").Value = "'  -0.22%  "
$ws.Range("D38").Value = "1.240.25"
$ws.Range("E38").Value = "'  -0.85%  "
$ws.Range("D39").Value = "'0.01801"
$ws.Range("E39").Value = "'  -0.23%  "
$ws.Range("E40").Value = "'  -1.18%  "
$ws.Range("D41").Value = "'6.247"
$ws.Range("E41").Value = "'  +2.53%  "
$ws.Range("D42").Value = "'0.8999"
$ws.Range("E42").Value = "'  -0.50%  "
$ws.Range("D43").Value = "'0.9990"
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("D44").Value = "'101.92"
$ws.Range("E44").Value = "'  +0.53%  "
$ws.Range("D45").Value = "'65.31"
$ws.Range("E45").Value = "'  -2.71%  "
$ws.Range("E46").Value = "'  +0.85%  "
$ws.Range("D47").Value = "'7.013"
$ws.Range("E47").Value = "'  -4.12%  "
$ws.Range("D48").Value = "'0.3993"
$ws.Range("E48").Value = "'  -0.38%  "
$ws.Range("D49").Value = "'1.667"
$ws.Range("E49").Value = "'  -1.66%  "
$ws.Range("D50").Value = "'8.891"
$ws.Range("E50").Value = "'  -2.93%  "
$ws.Range("D51").Value = "'0.1119"
$ws.Range("E51").Value = "'  -0.22%  "
